$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2:F61").Sort($ws.Range("E2:E61"))
